$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12 ---
$ws.Range("G12").Value = 2.25
$ws.Range("H12").Value = 2.88
$ws.Range("I12").Value = 3.6
$ws.Range("J12").Value = 3.1
$ws.Range("L12").Value = 4.33
$ws.Range("Z12").Value = 9.5
$ws.Range("AA12").Value = 10
$ws.Range("AB12").Value = 21
$ws.Range("AF12").Value = 6
$ws.Range("AH12").Value = 67
$ws.Range("AJ12").Value = 8
$ws.Range("AK12").Value = 17
$ws.Range("AL12").Value = 15
$ws.Range("AR12").Value = 2.03
$ws.Range("AS12").Value = 1.83

# --- Row 31 ---
$ws.Range("G31").Value = 2.25
$ws.Range("I31").Value = 3.3
$ws.Range("J31").Value = 3
$ws.Range("L31").Value = 4
$ws.Range("Q31").Value = 2.15
$ws.Range("R31").Value = 1.67
$ws.Range("W31").Value = 1.91
$ws.Range("X31").Value = 1.91
$ws.Range("Y31").Value = 7
$ws.Range("Z31").Value = 10
$ws.Range("AB31").Value = 21
$ws.Range("AC31").Value = 19
$ws.Range("AJ31").Value = 9.5
$ws.Range("AL31").Value = 12
$ws.Range("AN31").Value = 29
$ws.Range("AO31").Value = 41

# --- Row 32 ---
$ws.Range("G32").Value = 2.77
$ws.Range("H32").Value = 3.5
$ws.Range("I32").Value = 2.25
$ws.Range("J32").Value = 3.3
$ws.Range("K32").Value = 2.18
$ws.Range("L32").Value = 2.8
$ws.Range("S32").Value = 2.55
$ws.Range("W32").Value = 1.57
$ws.Range("X32").Value = 2.12
$ws.Range("Y32").Value = 10.75
$ws.Range("Z32").Value = 15.5
$ws.Range("AA32").Value = 10.25
$ws.Range("AB32").Value = 32
$ws.Range("AC32").Value = 22
$ws.Range("AD32").Value = 27
$ws.Range("AF32").Value = 6.9
$ws.Range("AG32").Value = 12.5
$ws.Range("AI32").Value = 300
$ws.Range("AJ32").Value = 9.5
$ws.Range("AK32").Value = 12
$ws.Range("AL32").Value = 9
$ws.Range("AM32").Value = 23
$ws.Range("AN32").Value = 17
$ws.Range("AO32").Value = 24
